# Generate Report for Handoff
# Updates the "b.md" row in each sheet to reflect that file "b" has now
# been handed off (status "Ready for handoff") with a new handoff file/
# timestamp, mirroring the existing "a" row's shape.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: row 3 is the b.md entry ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-27-19 02:27:26"

# ---- zh-cn sheet: row 3 is the b.md entry ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-19 02:27:24"

foreach ($h in $wsZh.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---- de-de sheet: row 3 is the b.md entry ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-19 02:27:26"

foreach ($h in $wsDe.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
